{"js": "// Replace the 100 multiplication-fact answers in the single table on the page.\n// Each entry is [rowIndex, colIndex, expectedOldText, newText]. We address cells\n// directly by (row, col) rather than searching for text, since several new\n// values collide with other cells' old values (e.g. \"31\u00d761=1891\" is both an old\n// value at one cell and a new value at another), which would make a global\n// text search-and-replace ambiguous/order-dependent. Using the paragraph's own\n// Range for insertText (rather than the cell body) preserves the existing run\n// (font/size) and paragraph (alignment) formatting.\n\nconst replacements = [\n  [0, 0, \"80\u00d740=3200\", \"59\u00d765=3835\"],\n  [0, 1, \"100\u00d743=4300\", \"86\u00d749=4214\"],\n  [0, 2, \"100\u00d721=2100\", \"48\u00d792=4416\"],\n  [0, 3, \"43\u00d757=2451\", \"93\u00d776=7068\"],\n  [0, 4, \"92\u00d711=1012\", \"100\u00d792=9200\"],\n  [1, 0, \"18\u00d780=1440\", \"13\u00d711=143\"],\n  [1, 1, \"78\u00d798=7644\", \"68\u00d723=1564\"],\n  [1, 2, \"37\u00d787=3219\", \"11\u00d710=110\"],\n  [1, 3, \"39\u00d734=1326\", \"33\u00d740=1320\"],\n  [1, 4, \"62\u00d760=3720\", \"10\u00d775=750\"],\n  [2, 0, \"48\u00d739=1872\", \"85\u00d734=2890\"],\n  [2, 1, \"91\u00d721=1911\", \"26\u00d731=806\"],\n  [2, 2, \"45\u00d746=2070\", \"45\u00d752=2340\"],\n  [2, 3, \"66\u00d726=1716\", \"28\u00d755=1540\"],\n  [2, 4, \"91\u00d714=1274\", \"83\u00d789=7387\"],\n  [3, 0, \"31\u00d761=1891\", \"84\u00d711=924\"],\n  [3, 1, \"12\u00d744=528\", \"68\u00d756=3808\"],\n  [3, 2, \"93\u00d753=4929\", \"12\u00d726=312\"],\n  [3, 3, \"31\u00d734=1054\", \"91\u00d794=8554\"],\n  [3, 4, \"54\u00d779=4266\", \"75\u00d724=1800\"],\n  [4, 0, \"74\u00d742=3108\", \"52\u00d770=3640\"],\n  [4, 1, \"48\u00d745=2160\", \"89\u00d755=4895\"],\n  [4, 2, \"13\u00d770=910\", \"85\u00d734=2890\"],\n  [4, 3, \"19\u00d785=1615\", \"20\u00d774=1480\"],\n  [4, 4, \"88\u00d784=7392\", \"41\u00d733=1353\"],\n  [5, 0, \"86\u00d770=6020\", \"46\u00d724=1104\"],\n  [5, 1, \"90\u00d751=4590\", \"46\u00d793=4278\"],\n  [5, 2, \"63\u00d723=1449\", \"56\u00d759=3304\"],\n  [5, 3, \"15\u00d747=705\", \"35\u00d7100=3500\"],\n  [5, 4, \"20\u00d747=940\", \"46\u00d716=736\"],\n  [6, 0, \"97\u00d713=1261\", \"37\u00d764=2368\"],\n  [6, 1, \"78\u00d7100=7800\", \"23\u00d726=598\"],\n  [6, 2, \"34\u00d796=3264\", \"28\u00d770=1960\"],\n  [6, 3, \"67\u00d721=1407\", \"72\u00d781=5832\"],\n  [6, 4, \"29\u00d723=667\", \"50\u00d753=2650\"],\n  [7, 0, \"17\u00d793=1581\", \"80\u00d770=5600\"],\n  [7, 1, \"79\u00d734=2686\", \"25\u00d737=925\"],\n  [7, 2, \"84\u00d775=6300\", \"98\u00d712=1176\"],\n  [7, 3, \"98\u00d790=8820\", \"39\u00d764=2496\"],\n  [7, 4, \"71\u00d712=852\", \"71\u00d773=5183\"],\n  [8, 0, \"92\u00d769=6348\", \"65\u00d729=1885\"],\n  [8, 1, \"50\u00d763=3150\", \"69\u00d796=6624\"],\n  [8, 2, \"74\u00d750=3700\", \"72\u00d766=4752\"],\n  [8, 3, \"71\u00d727=1917\", \"25\u00d733=825\"],\n  [8, 4, \"23\u00d792=2116\", \"55\u00d788=4840\"],\n  [9, 0, \"99\u00d790=8910\", \"91\u00d735=3185\"],\n  [9, 1, \"82\u00d753=4346\", \"76\u00d773=5548\"],\n  [9, 2, \"11\u00d727=297\", \"59\u00d725=1475\"],\n  [9, 3, \"24\u00d779=1896\", \"35\u00d751=1785\"],\n  [9, 4, \"25\u00d789=2225\", \"12\u00d777=924\"],\n  [10, 0, \"93\u00d769=6417\", \"79\u00d735=2765\"],\n  [10, 1, \"22\u00d740=880\", \"78\u00d771=5538\"],\n  [10, 2, \"23\u00d758=1334\", \"53\u00d794=4982\"],\n  [10, 3, \"45\u00d790=4050\", \"25\u00d797=2425\"],\n  [10, 4, \"37\u00d712=444\", \"61\u00d750=3050\"],\n  [11, 0, \"94\u00d764=6016\", \"97\u00d777=7469\"],\n  [11, 1, \"53\u00d754=2862\", \"60\u00d745=2700\"],\n  [11, 2, \"40\u00d743=1720\", \"40\u00d737=1480\"],\n  [11, 3, \"60\u00d772=4320\", \"44\u00d769=3036\"],\n  [11, 4, \"98\u00d767=6566\", \"31\u00d761=1891\"],\n  [12, 0, \"70\u00d734=2380\", \"25\u00d732=800\"],\n  [12, 1, \"25\u00d711=275\", \"36\u00d742=1512\"],\n  [12, 2, \"47\u00d760=2820\", \"62\u00d765=4030\"],\n  [12, 3, \"45\u00d714=630\", \"77\u00d769=5313\"],\n  [12, 4, \"14\u00d768=952\", \"58\u00d765=3770\"],\n  [13, 0, \"79\u00d724=1896\", \"81\u00d734=2754\"],\n  [13, 1, \"65\u00d754=3510\", \"85\u00d775=6375\"],\n  [13, 2, \"25\u00d775=1875\", \"30\u00d711=330\"],\n  [13, 3, \"86\u00d774=6364\", \"90\u00d723=2070\"],\n  [13, 4, \"37\u00d752=1924\", \"39\u00d722=858\"],\n  [14, 0, \"54\u00d762=3348\", \"53\u00d772=3816\"],\n  [14, 1, \"68\u00d734=2312\", \"64\u00d762=3968\"],\n  [14, 2, \"64\u00d728=1792\", \"24\u00d778=1872\"],\n  [14, 3, \"99\u00d730=2970\", \"34\u00d793=3162\"],\n  [14, 4, \"27\u00d794=2538\", \"10\u00d721=210\"],\n  [15, 0, \"20\u00d711=220\", \"79\u00d795=7505\"],\n  [15, 1, \"19\u00d712=228\", \"12\u00d783=996\"],\n  [15, 2, \"47\u00d777=3619\", \"50\u00d744=2200\"],\n  [15, 3, \"29\u00d775=2175\", \"49\u00d762=3038\"],\n  [15, 4, \"34\u00d774=2516\", \"93\u00d737=3441\"],\n  [16, 0, \"15\u00d785=1275\", \"47\u00d776=3572\"],\n  [16, 1, \"19\u00d731=589\", \"56\u00d771=3976\"],\n  [16, 2, \"35\u00d780=2800\", \"34\u00d777=2618\"],\n  [16, 3, \"82\u00d764=5248\", \"78\u00d773=5694\"],\n  [16, 4, \"33\u00d776=2508\", \"91\u00d784=7644\"],\n  [17, 0, \"42\u00d770=2940\", \"18\u00d736=648\"],\n  [17, 1, \"71\u00d779=5609\", \"61\u00d759=3599\"],\n  [17, 2, \"99\u00d788=8712\", \"14\u00d788=1232\"],\n  [17, 3, \"99\u00d792=9108\", \"25\u00d762=1550\"],\n  [17, 4, \"52\u00d798=5096\", \"26\u00d716=416\"],\n  [18, 0, \"96\u00d721=2016\", \"73\u00d799=7227\"],\n  [18, 1, \"83\u00d786=7138\", \"24\u00d754=1296\"],\n  [18, 2, \"47\u00d738=1786\", \"97\u00d756=5432\"],\n  [18, 3, \"65\u00d766=4290\", \"13\u00d767=871\"],\n  [18, 4, \"23\u00d773=1679\", \"27\u00d759=1593\"],\n  [19, 0, \"66\u00d736=2376\", \"46\u00d713=598\"],\n  [19, 1, \"15\u00d729=435\", \"91\u00d751=4641\"],\n  [19, 2, \"94\u00d726=2444\", \"83\u00d740=3320\"],\n  [19, 3, \"81\u00d763=5103\", \"14\u00d759=826\"],\n  [19, 4, \"29\u00d783=2407\", \"99\u00d746=4554\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\nconst table = tables.items[0];\n\n// Load all cells' paragraph ranges in one batch.\nconst cellInfo = replacements.map(([row, col]) => {\n  const cell = table.getCell(row, col);\n  const para = cell.body.paragraphs.getFirst();\n  para.load(\"text\");\n  const range = para.getRange();\n  return { para, range };\n});\n\nawait context.sync();\n\n// Verify the existing text matches what the diff expects before touching it,\n// then replace it in place so formatting (fonts, size, alignment) is kept.\nfor (let i = 0; i < replacements.length; i++) {\n  const [row, col, oldText, newText] = replacements[i];\n  const actual = cellInfo[i].para.text;\n  if (actual !== oldText) {\n    throw new Error(\n      `Cell (${row},${col}) expected \"${oldText}\" but found \"${actual}\"`\n    );\n  }\n  cellInfo[i].range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 multiplication-fact answers in the single table on the page.\n# Each entry addresses a cell directly by its (Row, Col) position (1-based, as\n# used by the Word object model) rather than searching for text. Several new\n# values collide with other cells' old values (e.g. \"31\u00d761=1891\" is both an\n# old value at one cell and a new value at another), which would make a\n# global Find/Replace pass ambiguous/order-dependent. Assigning directly to\n# Cell.Range.Text keeps the run/paragraph formatting (font, size, alignment)\n# already present in that cell intact.\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = '80\u00d740=3200'; NewText = '59\u00d765=3835' },\n    @{ Row = 1; Col = 2; OldText = '100\u00d743=4300'; NewText = '86\u00d749=4214' },\n    @{ Row = 1; Col = 3; OldText = '100\u00d721=2100'; NewText = '48\u00d792=4416' },\n    @{ Row = 1; Col = 4; OldText = '43\u00d757=2451'; NewText = '93\u00d776=7068' },\n    @{ Row = 1; Col = 5; OldText = '92\u00d711=1012'; NewText = '100\u00d792=9200' },\n    @{ Row = 2; Col = 1; OldText = '18\u00d780=1440'; NewText = '13\u00d711=143' },\n    @{ Row = 2; Col = 2; OldText = '78\u00d798=7644'; NewText = '68\u00d723=1564' },\n    @{ Row = 2; Col = 3; OldText = '37\u00d787=3219'; NewText = '11\u00d710=110' },\n    @{ Row = 2; Col = 4; OldText = '39\u00d734=1326'; NewText = '33\u00d740=1320' },\n    @{ Row = 2; Col = 5; OldText = '62\u00d760=3720'; NewText = '10\u00d775=750' },\n    @{ Row = 3; Col = 1; OldText = '48\u00d739=1872'; NewText = '85\u00d734=2890' },\n    @{ Row = 3; Col = 2; OldText = '91\u00d721=1911'; NewText = '26\u00d731=806' },\n    @{ Row = 3; Col = 3; OldText = '45\u00d746=2070'; NewText = '45\u00d752=2340' },\n    @{ Row = 3; Col = 4; OldText = '66\u00d726=1716'; NewText = '28\u00d755=1540' },\n    @{ Row = 3; Col = 5; OldText = '91\u00d714=1274'; NewText = '83\u00d789=7387' },\n    @{ Row = 4; Col = 1; OldText = '31\u00d761=1891'; NewText = '84\u00d711=924' },\n    @{ Row = 4; Col = 2; OldText = '12\u00d744=528'; NewText = '68\u00d756=3808' },\n    @{ Row = 4; Col = 3; OldText = '93\u00d753=4929'; NewText = '12\u00d726=312' },\n    @{ Row = 4; Col = 4; OldText = '31\u00d734=1054'; NewText = '91\u00d794=8554' },\n    @{ Row = 4; Col = 5; OldText = '54\u00d779=4266'; NewText = '75\u00d724=1800' },\n    @{ Row = 5; Col = 1; OldText = '74\u00d742=3108'; NewText = '52\u00d770=3640' },\n    @{ Row = 5; Col = 2; OldText = '48\u00d745=2160'; NewText = '89\u00d755=4895' },\n    @{ Row = 5; Col = 3; OldText = '13\u00d770=910'; NewText = '85\u00d734=2890' },\n    @{ Row = 5; Col = 4; OldText = '19\u00d785=1615'; NewText = '20\u00d774=1480' },\n    @{ Row = 5; Col = 5; OldText = '88\u00d784=7392'; NewText = '41\u00d733=1353' },\n    @{ Row = 6; Col = 1; OldText = '86\u00d770=6020'; NewText = '46\u00d724=1104' },\n    @{ Row = 6; Col = 2; OldText = '90\u00d751=4590'; NewText = '46\u00d793=4278' },\n    @{ Row = 6; Col = 3; OldText = '63\u00d723=1449'; NewText = '56\u00d759=3304' },\n    @{ Row = 6; Col = 4; OldText = '15\u00d747=705'; NewText = '35\u00d7100=3500' },\n    @{ Row = 6; Col = 5; OldText = '20\u00d747=940'; NewText = '46\u00d716=736' },\n    @{ Row = 7; Col = 1; OldText = '97\u00d713=1261'; NewText = '37\u00d764=2368' },\n    @{ Row = 7; Col = 2; OldText = '78\u00d7100=7800'; NewText = '23\u00d726=598' },\n    @{ Row = 7; Col = 3; OldText = '34\u00d796=3264'; NewText = '28\u00d770=1960' },\n    @{ Row = 7; Col = 4; OldText = '67\u00d721=1407'; NewText = '72\u00d781=5832' },\n    @{ Row = 7; Col = 5; OldText = '29\u00d723=667'; NewText = '50\u00d753=2650' },\n    @{ Row = 8; Col = 1; OldText = '17\u00d793=1581'; NewText = '80\u00d770=5600' },\n    @{ Row = 8; Col = 2; OldText = '79\u00d734=2686'; NewText = '25\u00d737=925' },\n    @{ Row = 8; Col = 3; OldText = '84\u00d775=6300'; NewText = '98\u00d712=1176' },\n    @{ Row = 8; Col = 4; OldText = '98\u00d790=8820'; NewText = '39\u00d764=2496' },\n    @{ Row = 8; Col = 5; OldText = '71\u00d712=852'; NewText = '71\u00d773=5183' },\n    @{ Row = 9; Col = 1; OldText = '92\u00d769=6348'; NewText = '65\u00d729=1885' },\n    @{ Row = 9; Col = 2; OldText = '50\u00d763=3150'; NewText = '69\u00d796=6624' },\n    @{ Row = 9; Col = 3; OldText = '74\u00d750=3700'; NewText = '72\u00d766=4752' },\n    @{ Row = 9; Col = 4; OldText = '71\u00d727=1917'; NewText = '25\u00d733=825' },\n    @{ Row = 9; Col = 5; OldText = '23\u00d792=2116'; NewText = '55\u00d788=4840' },\n    @{ Row = 10; Col = 1; OldText = '99\u00d790=8910'; NewText = '91\u00d735=3185' },\n    @{ Row = 10; Col = 2; OldText = '82\u00d753=4346'; NewText = '76\u00d773=5548' },\n    @{ Row = 10; Col = 3; OldText = '11\u00d727=297'; NewText = '59\u00d725=1475' },\n    @{ Row = 10; Col = 4; OldText = '24\u00d779=1896'; NewText = '35\u00d751=1785' },\n    @{ Row = 10; Col = 5; OldText = '25\u00d789=2225'; NewText = '12\u00d777=924' },\n    @{ Row = 11; Col = 1; OldText = '93\u00d769=6417'; NewText = '79\u00d735=2765' },\n    @{ Row = 11; Col = 2; OldText = '22\u00d740=880'; NewText = '78\u00d771=5538' },\n    @{ Row = 11; Col = 3; OldText = '23\u00d758=1334'; NewText = '53\u00d794=4982' },\n    @{ Row = 11; Col = 4; OldText = '45\u00d790=4050'; NewText = '25\u00d797=2425' },\n    @{ Row = 11; Col = 5; OldText = '37\u00d712=444'; NewText = '61\u00d750=3050' },\n    @{ Row = 12; Col = 1; OldText = '94\u00d764=6016'; NewText = '97\u00d777=7469' },\n    @{ Row = 12; Col = 2; OldText = '53\u00d754=2862'; NewText = '60\u00d745=2700' },\n    @{ Row = 12; Col = 3; OldText = '40\u00d743=1720'; NewText = '40\u00d737=1480' },\n    @{ Row = 12; Col = 4; OldText = '60\u00d772=4320'; NewText = '44\u00d769=3036' },\n    @{ Row = 12; Col = 5; OldText = '98\u00d767=6566'; NewText = '31\u00d761=1891' },\n    @{ Row = 13; Col = 1; OldText = '70\u00d734=2380'; NewText = '25\u00d732=800' },\n    @{ Row = 13; Col = 2; OldText = '25\u00d711=275'; NewText = '36\u00d742=1512' },\n    @{ Row = 13; Col = 3; OldText = '47\u00d760=2820'; NewText = '62\u00d765=4030' },\n    @{ Row = 13; Col = 4; OldText = '45\u00d714=630'; NewText = '77\u00d769=5313' },\n    @{ Row = 13; Col = 5; OldText = '14\u00d768=952'; NewText = '58\u00d765=3770' },\n    @{ Row = 14; Col = 1; OldText = '79\u00d724=1896'; NewText = '81\u00d734=2754' },\n    @{ Row = 14; Col = 2; OldText = '65\u00d754=3510'; NewText = '85\u00d775=6375' },\n    @{ Row = 14; Col = 3; OldText = '25\u00d775=1875'; NewText = '30\u00d711=330' },\n    @{ Row = 14; Col = 4; OldText = '86\u00d774=6364'; NewText = '90\u00d723=2070' },\n    @{ Row = 14; Col = 5; OldText = '37\u00d752=1924'; NewText = '39\u00d722=858' },\n    @{ Row = 15; Col = 1; OldText = '54\u00d762=3348'; NewText = '53\u00d772=3816' },\n    @{ Row = 15; Col = 2; OldText = '68\u00d734=2312'; NewText = '64\u00d762=3968' },\n    @{ Row = 15; Col = 3; OldText = '64\u00d728=1792'; NewText = '24\u00d778=1872' },\n    @{ Row = 15; Col = 4; OldText = '99\u00d730=2970'; NewText = '34\u00d793=3162' },\n    @{ Row = 15; Col = 5; OldText = '27\u00d794=2538'; NewText = '10\u00d721=210' },\n    @{ Row = 16; Col = 1; OldText = '20\u00d711=220'; NewText = '79\u00d795=7505' },\n    @{ Row = 16; Col = 2; OldText = '19\u00d712=228'; NewText = '12\u00d783=996' },\n    @{ Row = 16; Col = 3; OldText = '47\u00d777=3619'; NewText = '50\u00d744=2200' },\n    @{ Row = 16; Col = 4; OldText = '29\u00d775=2175'; NewText = '49\u00d762=3038' },\n    @{ Row = 16; Col = 5; OldText = '34\u00d774=2516'; NewText = '93\u00d737=3441' },\n    @{ Row = 17; Col = 1; OldText = '15\u00d785=1275'; NewText = '47\u00d776=3572' },\n    @{ Row = 17; Col = 2; OldText = '19\u00d731=589'; NewText = '56\u00d771=3976' },\n    @{ Row = 17; Col = 3; OldText = '35\u00d780=2800'; NewText = '34\u00d777=2618' },\n    @{ Row = 17; Col = 4; OldText = '82\u00d764=5248'; NewText = '78\u00d773=5694' },\n    @{ Row = 17; Col = 5; OldText = '33\u00d776=2508'; NewText = '91\u00d784=7644' },\n    @{ Row = 18; Col = 1; OldText = '42\u00d770=2940'; NewText = '18\u00d736=648' },\n    @{ Row = 18; Col = 2; OldText = '71\u00d779=5609'; NewText = '61\u00d759=3599' },\n    @{ Row = 18; Col = 3; OldText = '99\u00d788=8712'; NewText = '14\u00d788=1232' },\n    @{ Row = 18; Col = 4; OldText = '99\u00d792=9108'; NewText = '25\u00d762=1550' },\n    @{ Row = 18; Col = 5; OldText = '52\u00d798=5096'; NewText = '26\u00d716=416' },\n    @{ Row = 19; Col = 1; OldText = '96\u00d721=2016'; NewText = '73\u00d799=7227' },\n    @{ Row = 19; Col = 2; OldText = '83\u00d786=7138'; NewText = '24\u00d754=1296' },\n    @{ Row = 19; Col = 3; OldText = '47\u00d738=1786'; NewText = '97\u00d756=5432' },\n    @{ Row = 19; Col = 4; OldText = '65\u00d766=4290'; NewText = '13\u00d767=871' },\n    @{ Row = 19; Col = 5; OldText = '23\u00d773=1679'; NewText = '27\u00d759=1593' },\n    @{ Row = 20; Col = 1; OldText = '66\u00d736=2376'; NewText = '46\u00d713=598' },\n    @{ Row = 20; Col = 2; OldText = '15\u00d729=435'; NewText = '91\u00d751=4641' },\n    @{ Row = 20; Col = 3; OldText = '94\u00d726=2444'; NewText = '83\u00d740=3320' },\n    @{ Row = 20; Col = 4; OldText = '81\u00d763=5103'; NewText = '14\u00d759=826' },\n    @{ Row = 20; Col = 5; OldText = '29\u00d783=2407'; NewText = '99\u00d746=4554' }\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nforeach ($item in $replacements) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    $actual = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($actual -ne $item.OldText) {\n        throw \"Cell ($($item.Row),$($item.Col)) expected '$($item.OldText)' but found '$actual'\"\n    }\n    $cell.Range.Text = $item.NewText\n}\n"}
